$wb = $excel.ActiveWorkbook

# Add the new "InvalidLogin" worksheet after the existing "ValidLogin" sheet
$validSheet = $wb.Worksheets.Item("ValidLogin")
$newSheet = $wb.Worksheets.Add($null, $validSheet)
$newSheet.Name = "InvalidLogin"

# Header row (same headers as ValidLogin)
$newSheet.Range("A1").Value = "UserName"
$newSheet.Range("B1").Value = "Password"

# Invalid login test data
$newSheet.Range("A2").Value = "abcd"
$newSheet.Range("B2").Value = "xyz"

$newSheet.Range("B3").Select()

# Make InvalidLogin the active sheet/tab
$newSheet.Activate()
